$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 3546
$ws.Range("O3").Value = 4732
$ws.Range("O4").Value = 646
$ws.Range("O5").Value = 292
